# Implementing NLTK for portuguese text processing.
# The underlying "Query" sheet (the Power Query staging sheet holding the
# SQL used to pull the tweets) is no longer needed, and the results sheet
# is renamed to the more descriptive "model".

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

# Rename "Result 1" -> "model"
$ws1 = $wb.Worksheets.Item("Result 1")
$ws1.Name = "model"

# Drop the now-unused "Query" sheet entirely (also drops its only cell's
# shared-string reference, the SELECT ... FROM twitter.nome_social t text).
$ws2 = $wb.Worksheets.Item("Query")
$ws2.Delete()
